# S0155_DeleteCase.xlsx - replace the three case-id rows (A2:A4) with a new
# batch of case IDs, and update the sheet selection to reflect the rows that
# were just entered.
#
# The case IDs look like numbers ("00001323") but must stay text with their
# leading zeros intact, exactly like the original shared-string values did.
# A direct `Range("A2").Value = "00001323"` gets auto-coerced to the number
# 1323 (losing the leading zeros) unless the cell/column is already
# formatted as Text. So we stage each value in a scratch cell that has been
# explicitly formatted as Text, copy it, and paste it (values+formats) onto
# the destination cell - this is the standard Excel trick for typing
# leading-zero "numbers" that must remain text, and it keeps the written
# cell a plain shared-string entry (t="s") like the rest of the column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$scratch = $ws.Range("Z100")
$scratch.NumberFormat = "@"

$values = @("00001323", "00001326", "00001328")
$targets = @("A2", "A3", "A4")

for ($i = 0; $i -lt $values.Length; $i++) {
    $scratch.Value = $values[$i]
    $scratch.Copy()
    $ws.Range($targets[$i]).PasteSpecial(-4163)
}

$excel.CutCopyMode = $false

# Remove the scratch cell entirely so it doesn't linger in the used range.
$scratch.EntireRow.Delete()

# Match the new selection left behind in the saved file: A2:A4 selected,
# A2 active.
$ws.Range("A2:A4").Select()
$excel.ActiveCell = $ws.Range("A2")
